# Altered Template and placed other files needed
# Adds a new Q&A row (row 8) to the Errors template sheet:
#   A8 = error description, B8 = solution, C8 = guidance link
# This introduces 3 new shared strings, extends the used range to A1:C8,
# sets the new row's height, and moves the active selection the way it
# ended up after the edit (selection on C9, just past the new row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended at the bottom of the table.
$ws.Range("A8").Value = "Invalid block tag on line 10: 'static'. Did you forget to register or load this tag?"
$ws.Range("B8").Value = "Add {% load staticfiles %} at the top of the html page."
$ws.Range("C8").Value = "http://stackoverflow.com/questions/27886477/invalid-block-tag-static"

# Match the row height used for the new row.
$ws.Rows.Item(8).RowHeight = 101.45

# Leave the selection where it ended up after typing the new row.
$ws.Range("C9").Select()
